$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'302.42"
$ws.Range("E2").Value = "'2.00%"
$ws.Range("D3").Value = "'44.10"
$ws.Range("E3").Value = "'6.41%"
$ws.Range("D4").Value = "'5.080"
$ws.Range("E4").Value = "'0.82%"
$ws.Range("D5").Value = "'0.07702"
$ws.Range("E5").Value = "'3.16%"
$ws.Range("D6").Value = "'1.620"
$ws.Range("E6").Value = "'2.88%"
$ws.Range("E7").Value = "'12.62%"
$ws.Range("D8").Value = "'0.1278"
$ws.Range("E8").Value = "'7.72%"
$ws.Range("D9").Value = "'0.1865"
$ws.Range("E9").Value = "'2.14%"
$ws.Range("D10").Value = "'0.09236"
$ws.Range("E10").Value = "'4.21%"
$ws.Range("E11").Value = "'0.72%"
$ws.Range("D12").Value = "'0.1046"
$ws.Range("E12").Value = "'-0.43%"
$ws.Range("E13").Value = "'0.18%"
$ws.Range("D14").Value = "'0.005769"
$ws.Range("E14").Value = "'-2.77%"
$ws.Range("E15").Value = "'1,912.28%"
$ws.Range("D17").Value = "'4.420"
$ws.Range("E18").Value = "'-3.86%"
$ws.Range("D19").Value = "'0.3351"
$ws.Range("E19").Value = "'1.82%"
$ws.Range("D20").Value = "'8.650"
$ws.Range("E20").Value = "'9.76%"
$ws.Range("D21").Value = "'0.1400"
$ws.Range("E21").Value = "'-0.66%"
$ws.Range("D22").Value = "'0.3178"
$ws.Range("E22").Value = "'7.14%"
$ws.Range("D23").Value = "'0.04175"
$ws.Range("E23").Value = "'3.37%"
$ws.Range("E24").Value = "'1.59%"
$ws.Range("E25").Value = "'13.91%"
$ws.Range("D26").Value = "'0.0001349"
$ws.Range("E26").Value = "'9.73%"
$ws.Range("D38").Value = "'0.02485"
$ws.Range("E38").Value = "'3.86%"
$ws.Range("D39").Value = "'0.05301"
$ws.Range("E39").Value = "'1.80%"
$ws.Range("D40").Value = "'0.005930"
$ws.Range("E40").Value = "'-13.79%"
$ws.Range("D41").Value = "'0.007753"
$ws.Range("E41").Value = "'-0.44%"
$ws.Range("D42").Value = "'0.1347"
$ws.Range("E42").Value = "'2.01%"
$ws.Range("D43").Value = "'0.007339"
$ws.Range("E43").Value = "'-0.50%"
$ws.Range("D44").Value = "'0.007542"
$ws.Range("E44").Value = "'5.22%"
$ws.Range("D45").Value = "'0.3016"
$ws.Range("E45").Value = "'-6.21%"
$ws.Range("D46").Value = "'0.00006656"
$ws.Range("E46").Value = "'6.80%"
$ws.Range("E47").Value = "'-0.05%"
$ws.Range("E48").Value = "'-7.40%"
$ws.Range("E49").Value = "'0.02%"
$ws.Range("D50").Value = "'0.00002099"
$ws.Range("E50").Value = "'-0.05%"
$ws.Range("D51").Value = "'0.0001999"
$ws.Range("E51").Value = "'-0.05%"
